$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells that must remain text (numeric-looking price strings)
# are given an explicit Text number format before the value is set,
# then restored to General so Excel does not re-interpret them as numbers.

# --- Rows 2-46: price / volume(1h) updates ---
# Row 2
$ws.Range("D2").Value = "73.149.79"
$ws.Range("E2").Value = "  -0.06%  "

# Row 3
$ws.Range("D3").Value = "3.983.56"
$ws.Range("E3").Value = "  -1.58%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "612.16"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +4.76%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.93"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +13.38%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.693"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +0.32%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.803"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +5.87%  "

# Row 10
$ws.Range("E10").Value = "  +9.39%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.18"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +6.72%  "

# Row 12
$ws.Range("E12").Value = "  +4.14%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.80"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +6.41%  "

# Row 14
$ws.Range("D14").Value = "4.614.09"
$ws.Range("E14").Value = "  -1.63%  "

# Row 15
$ws.Range("D15").Value = "3.976.36"
$ws.Range("E15").Value = "  -1.63%  "

# Row 16
$ws.Range("E16").Value = "  +3.80%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.35"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +0.45%  "

# Row 18
$ws.Range("E18").Value = "  +1.21%  "

# Row 19
$ws.Range("D19").Value = "73.076.13"
$ws.Range("E19").Value = "  -0.11%  "

# Row 20
$ws.Range("E20").Value = "  -0.69%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "458.25"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +3.53%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.88"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +5.29%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "97.55"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +0.25%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.40"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -4.03%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.36"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -1.49%  "

# Row 26
$ws.Range("E26").Value = "  -1.43%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.41"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -1.36%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.68"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -2.70%  "

# Row 29
$ws.Range("E29").Value = "  -1.50%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.38"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -1.29%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.95"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +0.75%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "14.09"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +3.24%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "49.58"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +2.36%  "

# Row 34
$ws.Range("E34").Value = "  +16.90%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.130"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -3.12%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "69.53"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +2.86%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "634.36"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -8.00%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.432"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -3.33%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.45"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +2.61%  "

# Row 40
$ws.Range("E40").Value = "  -0.79%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +0.09%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.29"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +47.93%  "

# Row 43
$ws.Range("E43").Value = "  +0.03%  "

# Row 44
$ws.Range("E44").Value = "  -1.63%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.59"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -6.05%  "

# Row 46
$ws.Range("E46").Value = "  -0.33%  "

# --- Rows 47-51: full row rewrites (coin reordering / replacement) ---
# Row 47
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.01"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -10.03%  "

# Row 48
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000298"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +9.52%  "

# Row 49
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.67"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -3.26%  "

# Row 50
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.43"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +1.02%  "

# Row 51
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.01"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -1.77%  "

